$d = $word.ActiveDocument

# 1. "32. -Halim wanted to conduct:" -> "32... -Halim-wanted'to conduct:"
#    (use Find to locate, then set .Text directly so the straight apostrophe
#     is not auto-corrected into a curly/smart quote)
$rng = $d.Content
$rng.Find.Execute("32. -Halim wanted to conduct:", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$rng.Text = "32... -Halim-wanted'to conduct:"

# 2. Remove the paragraph containing the picture (2nd paragraph)
$d.Paragraphs(2).Range.Delete()

# 3. "(a) © Besides the size of the boxes," -> ". (a) + Besides the size of the boxes,"
$d.Content.Find.Execute("(a) © Besides the size of the boxes,", $true, $false, $false, $false, $false,
                         $true, 1, $false, ". (a) + Besides the size of the boxes,", 2)

# 4. "(b) The results of Halim's experiment showed that the plastic box remained" ->
#    "(b) The results of ralims experiment snowed that the piastic box remained"
$d.Content.Find.Execute("(b) The results of Halim" + [char]0x2019 + "s experiment showed that the plastic box remained",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false, "(b) The results of ralims experiment snowed that the piastic box remained", 2)

# 5. "Past lows arf goss wien beakaowm |_|" -> "|| eats room gse uhan ty break dons |"
$d.Content.Find.Execute("Past lows arf goss wien beakaowm |_|", $true, $false, $false, $false, $false,
                         $true, 1, $false, "|| eats room gse uhan ty break dons |", 2)

# 6. "Bho _. (c) State one way" -> "Bho _- (c) State one way"
$d.Content.Find.Execute("Bho _. (c) State one way", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Bho _- (c) State one way", 2)

# 7. ". environment. a . (1)" -> ". environment. a . (1]"
$d.Content.Find.Execute(". environment. a . (1)", $true, $false, $false, $false, $false,
                         $true, 1, $false, ". environment. a . (1]", 2)
